$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set standard_hours for employee 1 (Sky Lee) to 4 instead of 8
$ws.Range("D2").Value = 4

# Delete the entire row 5 (employee id 69 / "Thim"), shifting cells up
$ws.Rows("5:5").Delete()
